# Adds two new columns (boolean + boolean-as-string) between the existing
# "column_4_string" and "column_7_int_multiple" columns, shifting the old
# E/F/G columns to G/H/I; updates column widths, selection, dimension and
# page setup to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at E:F (old E,F,G shift right to G,H,I)
$ws.Range("E1:F1").EntireColumn.Insert()

# Headers for the two new columns
$ws.Range("E1").Value = "column_5_boolean"
$ws.Range("F1").Value = "column_6_boolean_as_string"

# column_5_boolean: real boolean values
$ws.Range("E2").Value = $true
$ws.Range("E3").Value = $false

# column_6_boolean_as_string: literal text "true" / "TRUE" (must NOT be
# auto-coerced to a Boolean cell) -- format as Text, then paste-special
# a computed string value so the literal characters are preserved verbatim.
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F3").NumberFormat = "@"

$helper = $ws.Range("Z1")

$helper.Formula = '=T("true")'
$helper.Copy()
$ws.Range("F2").PasteSpecial(-4163)

$helper.Formula = '=T("TRUE")'
$helper.Copy()
$ws.Range("F3").PasteSpecial(-4163)

$helper.Clear()

# Column widths (author's bestFit widths from the original Excel session)
$ws.Columns.Item(4).ColumnWidth = 13.998697916666666
$ws.Columns.Item(5).ColumnWidth = 15.998697916666666
$ws.Columns.Item(6).ColumnWidth = 24.666666666666668

# Selection moves to H1
[void]$ws.Range("H1").Select()

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$wb.Save()
